$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = '${skill}'
$ws.Range("B1").Value = '${resources}'
$ws.Range("C1").Value = '${duration}'
$ws.Range("D1").Value = '${batch}'
$ws.Range("E1").Value = '${exp}'
$ws.Range("F1").Value = '${remarks}'
$ws.Range("G1").Value = '${type}'

# ---- Row 2 ----
$ws.Range("A2").Value = "Fullstack"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "6 Months"
$ws.Range("D2").Value = "Batch A"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Ready to deploy"
$ws.Range("G2").Value = "valid"

# ---- Row 3 ----
$ws.Range("A3").Value = "Java"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "3 Months"
$ws.Range("D3").Value = "Batch B"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = "Missing resources"
$ws.Range("G3").Value = "invalid"

# ---- Row 4 ----
$ws.Range("A4").Value = "React"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "Batch C"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Missing duration"
$ws.Range("G4").Value = "invalid"

# ---- Row 5 ----
$ws.Range("A5").Value = "Python"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "4 Months"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "No batch specified"
$ws.Range("G5").Value = "invalid"

# ---- Styling ----
# Set black font color across the whole used range first (data rows style),
# then make the header row bold on top of that (header row style).
$allRange = $ws.Range("A1:G5")
$allRange.Font.Color = 0

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true

# ---- Selection ----
$ws.Range("I17").Select() | Out-Null
